# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 4
    4  = 8
    5  = 5
    6  = 4
    7  = 5
    8  = 6
    9  = 8
    10 = 8
    11 = 7
    12 = 5
    13 = 5
    14 = 4
    15 = 7
    16 = 6
    17 = 7
    18 = 5
    19 = 6
    20 = 6
    21 = 5
    22 = 7
    23 = 4
    24 = 5
    25 = 9
    26 = 4
    27 = 4
    28 = 3
    29 = 7
    30 = 7
    31 = 5
    32 = 6
    33 = 7
    34 = 4
    35 = 3
    36 = 3
    37 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
